$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.192.56'
$ws.Range('E2').Value = '  -1.89%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.041.44'
$ws.Range('E3').Value = '  -3.21%  '
$ws.Range('E4').Value = '  +0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.35'
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.660'
$ws.Range('E6').Value = '  +1.05%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.26'
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '63.70'
$ws.Range('E9').Value = '  +6.29%  '
$ws.Range('E10').Value = '  -1.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0743'
$ws.Range('E11').Value = '  -3.92%  '
$ws.Range('E12').Value = '  -3.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.900'
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.11'
$ws.Range('E14').Value = '  -7.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.338.59'
$ws.Range('E15').Value = '  -2.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.36'
$ws.Range('E16').Value = '  -4.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.043.49'
$ws.Range('E17').Value = '  -4.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.37'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '36.160.78'
$ws.Range('E19').Value = '  -1.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.19'
$ws.Range('E20').Value = '  -3.31%  '
$ws.Range('D21').Value = '0.0₃0852'
$ws.Range('E21').Value = '  -3.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '236.96'
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('E23').Value = '  -6.93%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.34'
$ws.Range('E25').Value = '  -3.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.25'
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.24'
$ws.Range('E27').Value = '  -6.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '164.17'
$ws.Range('E28').Value = '  -2.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.87'
$ws.Range('E29').Value = '  -6.35%  '
$ws.Range('E30').Value = '  -2.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.19'
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.95'
$ws.Range('E32').Value = '  -8.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0596'
$ws.Range('E33').Value = '  -2.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.40'
$ws.Range('E34').Value = '  -7.70%  '
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0868'
$ws.Range('E36').Value = '  +1.93%  '
$ws.Range('E37').Value = '  -1.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.20'
$ws.Range('E38').Value = '  -10.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.02'
$ws.Range('E39').Value = '  +1.23%  '
$ws.Range('E40').Value = '  -6.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.87'
$ws.Range('E41').Value = '  -1.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0214'
$ws.Range('E42').Value = '  -3.46%  '
$ws.Range('E43').Value = '  -7.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '93.05'
$ws.Range('E44').Value = '  -4.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0902'
$ws.Range('E45').Value = '  -6.61%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.364.36'
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.74'
$ws.Range('E47').Value = '  -3.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.38'
$ws.Range('E48').Value = '  +4.12%  '
$ws.Range('E49').Value = '  +0.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.27'
$ws.Range('E50').Value = '  -7.28%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.224.42'
$ws.Range('E51').Value = '  -2.77%  '

Write-Host "Applied changes"